# "change conduent to buck and edit plots"
#
# Sheet 1 ("Actuarial Firm Summary16_21"): rename the firm formerly split
# across "Scott Terando" / "Conduent (formerly Buck and/or Xerox)" to
# "Buck"; fix mis-labelled "Internal Actuarial Services" / "Perac" rows to
# their correct firm names; normalize "Usi Consulting Group" capitalization
# to "USI Consulting Group"; update the 2021 USI Consulting Group AAL/UAL
# figures; and drop the stray duplicate McGriff 2016 row (data shifts up).
#
# Sheet 2 ("Actuarial Firm Summary 3_2021"): re-sort the alphabetical
# top-firm pie-chart data now that "Conduent ..." became "Buck" and
# "Internal Actuarial Services"/"Perac" became "Public Employee Retirement
# Administration Commission".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet 1 : targeted row fixes ----------------------------------------

# Row 11: "Scott Terando" -> "CalPERS"
$ws1.Cells.Item(11, 1).Value = "CalPERS"

# Row 12: "Conduent (formerly Buck and/or Xerox)" -> "Buck"
$ws1.Cells.Item(12, 1).Value = "Buck"

# Row 16: "Internal Actuarial Services" -> "New York City Office of the Actuary"
$ws1.Cells.Item(16, 1).Value = "New York City Office of the Actuary"

# Row 18: "Conduent (formerly Buck and/or Xerox)" -> "Buck"
$ws1.Cells.Item(18, 1).Value = "Buck"

# Row 21: "Internal Actuarial Services" -> "Public Employee Retirement Administration Commission"
$ws1.Cells.Item(21, 1).Value = "Public Employee Retirement Administration Commission"

# Row 25: "Perac" -> "Public Employee Retirement Administration Commission"
$ws1.Cells.Item(25, 1).Value = "Public Employee Retirement Administration Commission"

# Row 32: USI Consulting Group 2021 -- label casing + refreshed figures
$ws1.Cells.Item(32, 1).Value = "USI Consulting Group"
$ws1.Cells.Item(32, 3).Value = 45311366578
$ws1.Cells.Item(32, 4).Value = 0.00415055355248314
$ws1.Cells.Item(32, 5).Value = 13595866937

# Row 62: "Usi Consulting Group" -> "USI Consulting Group" (casing fix only)
$ws1.Cells.Item(62, 1).Value = "USI Consulting Group"

# Drop the stray duplicate "Usi Consulting Group" 2021 row (old row 65); the
# Principal Financial Group / McGriff 2016 rows below it shift up to 65/66.
$ws1.Rows.Item(65).Delete()

# ---- Sheet 2 : rebuild the (now re-alphabetized) summary table ----------

$sheet2Data = @(
    @("Buck", 161404206940, 0.0267678729795161, 45675353140, 2021),
    @("CalPERS", 587976000000, 0.0975121230195361, 110653248000, 2021),
    @("Cavanaugh Macdonald Consulting", 901490366670, 0.149506509680075, 154718673140, 2021),
    @("Cheiron", 357816501213, 0.0593416171488327, 97309179763, 2021),
    @("Foster & Foster", 98364489600, 0.0163131321867373, 28795731600, 2021),
    @("Gabriel, Roeder, Smith & Company (GRS)", 1448217745970, 0.240177808284865, 261234414860, 2021),
    @("Korn Ferry Hay Group", 52926848000, 0.00877758499192533, 12695832000, 2021),
    @("Milliman", 791425215550, 0.131252896341814, 49509898210, 2021),
    @("New York City Office of the Actuary", 248994048500, 0.04129409752858, 7493808500, 2021),
    @("New York State and Local Retirement Systems' Actuary", 231904000000, 0.0384598204292735, -28177076000, 2021),
    @("Nystrs Office Of The Actuary", 131077400000, 0.021738362711881, -17071064000, 2021),
    @("Office of The State Actuary - Washington", 97128000000, 0.0161080681603356, -13019696000, 2021),
    @("Public Employee Retirement Administration Commission", 104534296000, 0.0173363557888631, 31406700000, 2021),
    @("Pwc", 59395150000, 0.00985031221268181, 2321468000, 2021),
    @("Segal", 649806185770, 0.107766270605705, 160198525120, 2021),
    @("Others", 107312888761, 0.0177971679293788, 15051100475, 2021)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $r = $i + 2
    $row = $sheet2Data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
}
